# Nieuwe data toegevoegd via Streamlit op 2024-12-03 18:02:05
# Append one new record to the bottom of Sheet1's data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row below the existing data (column A is always
# populated for real records), then write the new record there.
$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$ws.Range("A$row").Value = "CompaNanny"
$ws.Range("B$row").Value = "CompaNanny Statenkwartier BSO"
$ws.Range("C$row").Value = "VGO"

# Column D holds free-form date-like text (mixed formats already present,
# e.g. "2024-09-02" vs "2023-08-23 00:00:00") rather than real Excel dates,
# so force text formatting first to stop auto-coercion into a date serial.
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "2024-06-24"

$ws.Range("E$row").Value = 0
$ws.Range("F$row").Value = 0
$ws.Range("G$row").Value = 1
$ws.Range("H$row").Value = 0
$ws.Range("I$row").Value = 0
$ws.Range("J$row").Value = 0
